# The deck's slide master currently uses the "Integral" theme (Red Violet
# color scheme, ppt/theme/theme1.xml). The target edit swaps the theme
# colors so the (only reachable/visible) master theme takes on the
# default "Office Theme" color palette -- dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink -- matching what used to live in ppt/theme/theme2.xml
# (which is only wired to the notes master and isn't independently
# reachable through the PowerPoint object model).
#
# ThemeColorScheme.Colors(index) follows DrawingML <clrScheme> order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# RGB is written using the standard VB/COM packed-BGR integer
# (R + G*256 + B*65536), so we convert each target hex color below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeHexColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeHexColors.Count; $i++) {
    $hex = $officeHexColors[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgrInt = $r + ($g * 256) + ($b * 65536)

    $color = $tcs.Colors($i + 1)
    $color.RGB = $bgrInt
}
